$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(3)

# Reposition / resize the oval ("a bola laranja")
$shape.Left = 643.5692125984252
$shape.Top = 294.64614173228347
$shape.Width = 167.2615748031496
$shape.Height = 165.04614173228347

# Give it a solid orange fill using the theme's accent2 scheme color
$shape.Fill.Visible = $true
$shape.Fill.Solid()
$shape.Fill.ForeColor.ObjectThemeColor = 6
